$wb = $excel.ActiveWorkbook

# --- Sheet "Matriz" ---
$ws1 = $wb.Worksheets.Item("Matriz")

# Header row
$ws1.Range("C1").Value = "X3"

# Data rows 2-5 updated values
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 1

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 0

$ws1.Range("A4").Value = 0
$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = 0

$ws1.Range("A5").Value = 0
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 0

# Remove rows 6 and 7 (old extra data), shrinking the used range to A1:C5
$ws1.Rows("6:7").Delete()

# --- Sheet "Pesos" ---
$ws2 = $wb.Worksheets.Item("Pesos")

$ws2.Range("A2").Value = 0.7226984551978297
$ws2.Range("A3").Value = 0.8104254099121104
$ws2.Range("A4").Value = -0.7706304189360769

# --- Sheet "Umbrales" ---
$ws3 = $wb.Worksheets.Item("Umbrales")

$ws3.Range("A2").Value = 0.6694997070989952

# --- Sheet "Configuracion" ---
$ws4 = $wb.Worksheets.Item("Configuracion")

$ws4.Range("A2").Value = "ESCALON"
